$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.148.12"
$ws.Range("E2").Value = "  +2.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.821.07"
$ws.Range("E3").Value = "  +0.71%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "630.28"
$ws.Range("E5").Value = "  +5.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.12"
$ws.Range("E6").Value = "  +0.52%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.819.38"
$ws.Range("E7").Value = "  +0.77%  "

$ws.Range("E8").Value = "  -0.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  +1.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.455"
$ws.Range("E11").Value = "  +0.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.62"
$ws.Range("E12").Value = "  +2.73%  "

$ws.Range("E13").Value = "  +0.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.06"
$ws.Range("E14").Value = "  +0.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.459.01"
$ws.Range("E15").Value = "  +0.64%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.708.25"
$ws.Range("E16").Value = "  -2.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.114.14"
$ws.Range("E17").Value = "  +1.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.15"
$ws.Range("E18").Value = "  -1.43%  "

$ws.Range("E19").Value = "  +1.06%  "

$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.60"
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.68"
$ws.Range("E22").Value = "  -1.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.710"
$ws.Range("E23").Value = "  +1.29%  "

$ws.Range("E24").Value = "  +4.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.87"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.03"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.16"
$ws.Range("E27").Value = "  +2.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("E28").Value = "  +0.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.968.48"
$ws.Range("E30").Value = "  +0.67%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.68"
$ws.Range("E31").Value = "  -2.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.24"
$ws.Range("E32").Value = "  +1.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.32"
$ws.Range("E33").Value = "  -1.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.25"
$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.10"
$ws.Range("E36").Value = "  +0.52%  "

$ws.Range("E37").Value = "  +2.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.149"
$ws.Range("E38").Value = "  +7.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.43"
$ws.Range("E39").Value = "  +5.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.94"
$ws.Range("E40").Value = "  +3.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.980"
$ws.Range("E41").Value = "  -0.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.47"
$ws.Range("E44").Value = "  +6.35%  "

$ws.Range("E45").Value = "  +0.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "154.45"
$ws.Range("E46").Value = "  +2.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.92"
$ws.Range("E47").Value = "  -1.52%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "42.67"
$ws.Range("E48").Value = "  -4.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.46"
$ws.Range("E49").Value = "  +1.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.89"
$ws.Range("E50").Value = "  +2.16%  "

$ws.Range("E51").Value = "  +8.79%  "
